$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDBbind v2013 benchmark (N=592)")

$ws.Range("C2").Value = 1.75
$ws.Range("D2").Value = 1.63
$ws.Range("E2").Value = 0.579
$ws.Range("F2").Value = 0.592

$ws.Range("C3").Value = 1.62
$ws.Range("D3").Value = 1.62
$ws.Range("E3").Value = 0.588
$ws.Range("F3").Value = 0.598

$ws.Range("C4").Value = 1.62
$ws.Range("D4").Value = 1.62
$ws.Range("E4").Value = 0.588
$ws.Range("F4").Value = 0.598

$ws.Range("C5").Value = 1.61
$ws.Range("D5").Value = 1.61
$ws.Range("E5").Value = 0.59
$ws.Range("F5").Value = 0.6

$ws.Range("C6").Value = 1.62
$ws.Range("D6").Value = 1.62
$ws.Range("E6").Value = 0.589
$ws.Range("F6").Value = 0.6

$ws.Range("C7").Value = 1.59
$ws.Range("D7").Value = 1.59
$ws.Range("E7").Value = 0.609
$ws.Range("F7").Value = 0.6

$ws.Range("C8").Value = 1.56
$ws.Range("D8").Value = 1.55
$ws.Range("E8").Value = 0.629
$ws.Range("F8").Value = 0.619

$ws.Range("C9").Value = 1.54
$ws.Range("D9").Value = 1.53
$ws.Range("E9").Value = 0.642
$ws.Range("F9").Value = 0.628

$ws.Range("C10").Value = 1.51
$ws.Range("D10").Value = 1.51
$ws.Range("E10").Value = 0.657
$ws.Range("F10").Value = 0.641

$ws.Range("C11").Value = 1.49
$ws.Range("D11").Value = 1.49
$ws.Range("E11").Value = 0.669
$ws.Range("F11").Value = 0.663

$ws.Range("C12").Value = 1.43
$ws.Range("D12").Value = 1.42
$ws.Range("E12").Value = 0.702
$ws.Range("F12").Value = 0.698

$ws.Range("C13").Value = 1.41
$ws.Range("D13").Value = 1.4
$ws.Range("E13").Value = 0.716
$ws.Range("F13").Value = 0.71

$ws.Range("C14").Value = 1.38
$ws.Range("D14").Value = 1.37
$ws.Range("E14").Value = 0.73
$ws.Range("F14").Value = 0.725

$ws.Range("C15").Value = 1.42
$ws.Range("D15").Value = 1.41
$ws.Range("E15").Value = 0.71
$ws.Range("F15").Value = 0.707

$ws.Range("C16").Value = 1.36
$ws.Range("D16").Value = 1.35
$ws.Range("E16").Value = 0.737
$ws.Range("F16").Value = 0.739

$ws.Range("C17").Value = 1.32
$ws.Range("D17").Value = 1.31
$ws.Range("E17").Value = 0.755
$ws.Range("F17").Value = 0.76

$ws.Range("C18").Value = 1.3
$ws.Range("D18").Value = 1.29
$ws.Range("E18").Value = 0.764
$ws.Range("F18").Value = 0.766
